$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 78125840
$ws.Range("I98").Value = 96154770
$ws.Range("J98").Value = 500
$ws.Range("K98").Value = 96154770
$ws.Range("L98").Value = 500
$ws.Range("M98").Value = -96153272
$ws.Range("N98").Value = -3496
$ws.Range("H113").Value = 3665644.5
$ws.Range("I113").Value = 11906986
$ws.Range("K113").Value = 11906986
$ws.Range("M113").Value = -11903732
$ws.Range("H122").Value = 78125840
$ws.Range("I122").Value = 96154770
$ws.Range("J122").Value = 500
$ws.Range("K122").Value = 288464310
$ws.Range("L122").Value = 1500
$ws.Range("M122").Value = -288461860
$ws.Range("N122").Value = -6400
$ws.Range("H132").Value = 24198462
$ws.Range("I132").Value = 26791084
$ws.Range("K132").Value = 80373252
$ws.Range("M132").Value = -80370722
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5561038.5
$ws.Range("I32").Value = 4455.7017
$ws.Range("K32").Value = 4455.7017
$ws.Range("M32").Value = -4168.7017
$ws.Range("H110").Value = 1797.7693
$ws.Range("I110").Value = 1328.8889
$ws.Range("J110").Value = 2852.75
$ws.Range("K110").Value = 1328.8889
$ws.Range("L110").Value = 2852.75
$ws.Range("M110").Value = 716.1111000000001
$ws.Range("N110").Value = -6942.75
$ws.Range("H123").Value = 30429
$ws.Range("J123").Value = 30429
$ws.Range("L123").Value = 30429
$ws.Range("N123").Value = -40229
$ws.Range("H132").Value = 906495.3
$ws.Range("I132").Value = 1417.7858
$ws.Range("J132").Value = 6538088.5
$ws.Range("K132").Value = 4253.357400000001
$ws.Range("L132").Value = 19614265.5
$ws.Range("M132").Value = -1723.357400000001
$ws.Range("N132").Value = -19619325.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 940
$ws.Range("I29").Value = 940
$ws.Range("K29").Value = 940
$ws.Range("M29").Value = -651
$ws.Range("H105").Value = 983.83
$ws.Range("I105").Value = 985.13513
$ws.Range("J105").Value = 980.11536
$ws.Range("K105").Value = 985.13513
$ws.Range("L105").Value = 980.11536
$ws.Range("M105").Value = 761.86487
$ws.Range("N105").Value = -4474.11536
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1288.3231
$ws.Range("I31").Value = 671.17145
$ws.Range("J31").Value = 2008.3334
$ws.Range("K31").Value = 671.17145
$ws.Range("L31").Value = 2008.3334
$ws.Range("M31").Value = -376.17145
$ws.Range("N31").Value = -2598.3334
$ws.Range("H34").Value = 1288.3231
$ws.Range("I34").Value = 671.17145
$ws.Range("J34").Value = 2008.3334
$ws.Range("K34").Value = 671.17145
$ws.Range("L34").Value = 2008.3334
$ws.Range("M34").Value = -469.17145
$ws.Range("N34").Value = -2412.3334
$ws.Range("H58").Value = 38462692
$ws.Range("I58").Value = 50000904
$ws.Range("J58").Value = 1983.3334
$ws.Range("K58").Value = 50000904
$ws.Range("L58").Value = 1983.3334
$ws.Range("M58").Value = -50000701
$ws.Range("N58").Value = -2389.3334
$ws.Range("H70").Value = 31578.6
$ws.Range("J70").Value = 31578.6
$ws.Range("L70").Value = 31578.6
$ws.Range("N70").Value = -32208.6
$ws.Range("H73").Value = 31578.6
$ws.Range("J73").Value = 31578.6
$ws.Range("L73").Value = 31578.6
$ws.Range("N73").Value = -33762.6
$ws.Range("H86").Value = 33369846
$ws.Range("I86").Value = 62502860
$ws.Range("J86").Value = 74972.57000000001
$ws.Range("K86").Value = 62502860
$ws.Range("L86").Value = 74972.57000000001
$ws.Range("M86").Value = -62501737
$ws.Range("N86").Value = -77218.57000000001
$ws.Range("H89").Value = 33369846
$ws.Range("I89").Value = 62502860
$ws.Range("J89").Value = 74972.57000000001
$ws.Range("K89").Value = 312514300
$ws.Range("L89").Value = 374862.85
$ws.Range("M89").Value = -312508684
$ws.Range("N89").Value = -386094.85
$ws.Range("H99").Value = 38467996
$ws.Range("I99").Value = 125014290
$ws.Range("K99").Value = 125014290
$ws.Range("M99").Value = -125012792
$ws.Range("H126").Value = 38467996
$ws.Range("I126").Value = 125014290
$ws.Range("K126").Value = 375042870
$ws.Range("M126").Value = -375040400
$ws.Range("H136").Value = 38462692
$ws.Range("I136").Value = 50000904
$ws.Range("J136").Value = 1983.3334
$ws.Range("K136").Value = 150002712
$ws.Range("L136").Value = 5950.0002
$ws.Range("M136").Value = -150000162
$ws.Range("N136").Value = -11050.0002
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 74.57143000000001
$ws.Range("I7").Value = 32.5
$ws.Range("J7").Value = 130.66667
$ws.Range("K7").Value = 97.5
$ws.Range("L7").Value = 392.00001
$ws.Range("M7").Value = 14.5
$ws.Range("N7").Value = -616.00001
$ws.Range("H107").Value = 2664
$ws.Range("I107").Value = 214.11765
$ws.Range("J107").Value = 3184.6
$ws.Range("K107").Value = 642.35295
$ws.Range("L107").Value = 9553.799999999999
$ws.Range("M107").Value = 1277.64705
$ws.Range("N107").Value = -13393.8
$ws.Range("H131").Value = 1838.7
$ws.Range("I131").Value = 3537.2632
$ws.Range("J131").Value = 1440.2716
$ws.Range("K131").Value = 10611.7896
$ws.Range("L131").Value = 4320.8148
$ws.Range("M131").Value = -5571.7896
$ws.Range("N131").Value = -14400.8148
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").ClearContents()
$ws.Range("N20").Value = 0
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("H43").Value = 7000
$ws.Range("I43").Value = 7000
$ws.Range("K43").Value = 7000
$ws.Range("M43").Value = -6849
$ws.Range("H57").Value = 14763.9375
$ws.Range("J57").Value = 14763.9375
$ws.Range("L57").Value = 14763.9375
$ws.Range("N57").Value = -16403.9375
$ws.Range("H80").Value = 1733
$ws.Range("I80").Value = 1679.8
$ws.Range("J80").Value = 1999
$ws.Range("K80").Value = 1679.8
$ws.Range("L80").Value = 1999
$ws.Range("M80").Value = -681.8
$ws.Range("N80").Value = -3995
$ws.Range("H83").Value = 1733
$ws.Range("I83").Value = 1679.8
$ws.Range("J83").Value = 1999
$ws.Range("K83").Value = 8399
$ws.Range("L83").Value = 9995
$ws.Range("M83").Value = -3407
$ws.Range("N83").Value = -19979
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 879.8
$ws.Range("I16").Value = 800
$ws.Range("J16").Value = 999.5
$ws.Range("K16").Value = 800
$ws.Range("L16").Value = 999.5
$ws.Range("M16").Value = -630
$ws.Range("N16").Value = -1339.5
$ws.Range("H31").Value = 4458.5
$ws.Range("I31").Value = 900
$ws.Range("J31").Value = 4853.8887
$ws.Range("K31").Value = 900
$ws.Range("L31").Value = 4853.8887
$ws.Range("M31").Value = -652
$ws.Range("N31").Value = -5349.8887
$ws.Range("H122").Value = 21780
$ws.Range("I122").Value = 26475
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 79425
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -76975
$ws.Range("N122").Value = -13900
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2039.3334
$ws.Range("I96").Value = 1547.2
$ws.Range("J96").Value = 4500
$ws.Range("K96").Value = 1547.2
$ws.Range("L96").Value = 4500
$ws.Range("M96").Value = -174.2
$ws.Range("N96").Value = -7246
$ws.Range("H107").Value = 488.2143
$ws.Range("I107").Value = 319.66666
$ws.Range("J107").Value = 1499.5
$ws.Range("K107").Value = 958.9999799999999
$ws.Range("L107").Value = 4498.5
$ws.Range("M107").Value = 961.0000200000001
$ws.Range("N107").Value = -8338.5
$ws.Range("H122").Value = 20244.186
$ws.Range("I122").Value = 24675.857
$ws.Range("J122").Value = 4733.3335
$ws.Range("K122").Value = 74027.571
$ws.Range("L122").Value = 14200.0005
$ws.Range("M122").Value = -71577.571
$ws.Range("N122").Value = -19100.0005
$ws.Range("H123").Value = 20000
$ws.Range("J123").Value = 20000
$ws.Range("L123").Value = 20000
$ws.Range("N123").Value = -29800
$ws.Range("H136").Value = 23812780
$ws.Range("I136").Value = 62502424
$ws.Range("J136").Value = 3767.6924
$ws.Range("K136").Value = 187507272
$ws.Range("L136").Value = 11303.0772
$ws.Range("M136").Value = -187504722
$ws.Range("N136").Value = -16403.0772
$ws.Range("H138").Value = 48464.5
$ws.Range("J138").Value = 48464.5
$ws.Range("L138").Value = 48464.5
$ws.Range("N138").Value = -58744.5
